$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Shared string label update: 'Copper ores and concentrates' -> 'Copper'
$ws1.Range("C4").Value = "Copper"
$ws2.Range("C4").Value = "Copper"

# Sheet 1, row 2
$row1_2 = @(261.9472600640528,302.5639882458125,500.2023325229081,370.5884416084285,948.980296248033,545.2749069198272,608.4981381272771,699.9134658860155,557.8761921188465,684.5176110887784,345.4458354396534,383.592260127869,770.9808728758536,778.0794695170396,785.3070845409572,1971.518226153727,1975.059190640531,1980.464632210437,1988.12750730279,1996.071621544525,1831.945142313142,1841.224274209226,1849.453327985607,1858.32322180638,1874.007768512969,549.8257572310469,574.458845297152,598.6690553282052,603.0295774294713,602.4273734116684,421.8495022842388,424.9260422441258,427.3967500290863,432.4065427301878,457.1720689401727,501.6851140183748,550.1633927447689,614.9605418309601,695.0946438025013,741.2807388388734,755.0702307087299,764.0366403115574,768.1352492600954,767.7590041496562,779.2561756643512,781.3430667654021,980.7295579684504,936.9579014828611,875.4125073453,849.7328494058657,827.9111981306293,826.0890175584319,816.5733055982598,805.6351036371897,801.7658895454782,780.9548675548587,796.3316069689852,826.0402590930486,869.7709067303512,923.1575793743789,962.3554788317695,1017.066885784799,1054.017227239844,1072.278744278378,1081.570337391527,1074.334070674188,1079.483747764993,1084.165350081036,1086.572255494762,1079.73414124064,1042.111477830739,1012.647486827561,997.4011763113591,995.2876280286404,994.9618308900295,978.2763238704088,976.2500150019441,979.2662918949696,990.5725417712317,1013.000167816229,1032.008258268781,1076.795214718335,1126.857748734077,1176.157693580835,1218.333346102828,1226.889707917357,1244.863300874676,1255.002920914531,1261.575287593672,1266.612866201052)
for ($c = 4; $c -le 93; $c++) { $ws1.Cells.Item(2, $c).Value = $row1_2[$c - 4] }

# Sheet 1, row 3
$row1_3 = @(16.48537993724342,19.15588856276863,32.31100287195206,23.60655368172608,62.55259237644246,35.4266080237992,39.56040732986385,45.54701273377307,36.05219760501377,44.48477195818336,21.61883087004233,24.31857910209132,50.25329065818438,50.69382503739336,51.14318108949972,130.7435212352095,130.9363135840933,131.2542657743722,131.7238117521813,132.2122502263641,121.1478155400725,121.7275831153449,122.2368904470552,122.7892485067165,123.7992289571712,34.84669890175937,36.45653176670924,38.03791239260138,38.28618390332246,38.20111887702271,26.03162858441701,26.19176717477905,26.31085366282995,26.59994426740899,28.21515432041273,31.15564098198947,34.36215401238877,38.66389717068211,43.99504061930835,47.04771391204812,47.92453729165306,48.47789529750027,48.70483968609841,48.63178038159134,49.35393159843112,49.44202486963682,62.78181388152785,59.79340307780486,55.61235404829592,53.83860428514903,52.31824961893121,52.14051987478872,51.44609530992616,50.65593303188957,50.33984841549508,48.88163758781356,49.85267362901053,51.78556820578853,54.65958901487823,58.18162948287311,60.74749083466683,64.35503583934266,66.77048649412571,67.93162187280578,68.4908011344737,67.93552344254115,68.2120638127667,68.45737380255733,68.55023677296357,68.02284589525068,65.42487082173695,63.37498644751282,62.27939032120079,62.06514445741061,61.97074115350754,60.77337882582801,60.56014968011761,60.68523264032837,61.36659723143823,62.79434806849929,63.98718790147537,66.91066426310346,70.18818927256976,73.41446114700635,76.16250020168715,76.64884506215728,77.76779697555394,78.36097092525588,78.71477060629303,78.96564121194649)
for ($c = 4; $c -le 93; $c++) { $ws1.Cells.Item(3, $c).Value = $row1_3[$c - 4] }

# Sheet 1, row 4
$row1_4 = @(804639.2050553659,548723.3242650117,355073.8984457181,291545.9398878629,328006.5025128267,256108.6506172263,249600.0595948292,237583.1566489558,305969.2535715514,280503.061528866,306576.0581405597,386749.5594112837,1801530.329054188,1808597.852628684,1812860.879121145,5689217.357277946,5691608.613121784,5695171.818355992,5699108.583822033,5703040.633327082,5374888.169368589,5379875.400106831,5383758.960864661,5387498.581991775,5392277.297044988,897654.9198817867,904451.2418461674,912703.4926970723,921745.5918088916,932509.0391692072,150868.9783683423,172187.626735935,201942.4464327527,243124.6617538523,298578.4961583769,365844.614574245,439540.0530973302,527759.1127925945,627695.6467452359,699670.2434188125,742043.9918424489,773660.4427390648,792760.8092025765,798438.0008717187,1237289.780110251,1345377.777379535,1480175.270517842,1407677.708723225,1313767.133457336,1245326.215671117,1118797.163870912,1079580.730513598,1041969.167460487,1010468.945903504,994277.2102867065,911996.6945050539,929012.1876663502,961220.8135492017,1006686.222834959,1062819.259188419,1063194.262851791,1124871.042273574,1173161.728744209,1205994.011840343,1229635.212015672,1203738.701838366,1216565.403806582,1224074.295360831,1223971.940356318,1210287.384748435,1133946.08153552,1091815.487368502,1058835.384802549,1036688.000889484,1017939.06019074,954355.7862108706,938489.0124144518,929601.949079862,930351.3488015354,941965.1663468548,913403.3500767406,944376.6032661406,981437.8623258495,1020354.627749699,1055891.942685894,1013997.38214646,1032381.711033423,1044412.968818649,1052575.503680638,1057647.212430665)
for ($c = 4; $c -le 93; $c++) { $ws1.Cells.Item(4, $c).Value = $row1_4[$c - 4] }

# Sheet 1, row 5
$row1_5 = @(1694765.500019309,995923.9418223609,499080.1759011137,306563.2458506891,287214.6565754326,203379.0319039942,165862.260963542,213790.6907842293,388716.5891332898,338415.2376235322,425186.3321593324,533217.8721845702,3529921.410042091,3532417.325437126,3535471.199540295,11253285.54521802,11257574.47100712,11262471.61563974,11267936.7948878,11273906.90943887,10772414.51871099,10779203.12676499,10786371.77143875,10794029.69298962,10802458.24477617,1801254.759697539,1813095.067817081,1828204.317212147,1848052.522558907,1874344.534818555,204374.3400486798,248792.9832535175,304400.5749421639,371886.2405280492,451115.5566156618,540961.4961567003,639195.4324622758,742479.5282097947,846440.2564803202,945879.6276757841,1035187.589359238,1108831.44948511,1161968.350506955,1191048.13782815,2343871.558658456,2646835.624409389,2601279.903456337,2536513.620622267,2457951.176699388,2371809.991789562,2117548.629441518,2035494.60977754,1964248.707069941,1908395.115549221,1871095.047059432,1689075.50539843,1691652.829119179,1712341.925468352,1748347.685928397,1796275.679387742,1728646.507207921,1789789.872622165,1852525.240778111,1913850.429029914,1971059.617246651,1935377.18952197,1977665.333512838,2010002.21939749,2031398.300048313,2041419.547921177,1955903.588918315,1944285.024792138,1923634.813189436,1895923.32212873,1863592.090538655,1740872.285278122,1707663.575576252,1678114.746642368,1654506.100357424,1638527.891577693,1531345.607897239,1532903.592032054,1542971.620552969,1560602.575033618,1584421.435270553,1476322.580884,1507393.412089416,1539339.19854678,1570411.786038188,1599045.674785746)
for ($c = 4; $c -le 93; $c++) { $ws1.Cells.Item(5, $c).Value = $row1_5[$c - 4] }

# Sheet 2, row 2
$row2_2 = @(261.9472600640528,564.5112483098653,1064.713580832773,1435.302022441202,2384.282318689235,2929.557225609062,3538.055363736339,4237.968829622355,4795.845021741201,5480.36263282998,5825.808468269634,6209.400728397502,6980.381601273356,7758.461070790396,8543.768155331352,10515.28638148508,12490.34557212561,14470.81020433605,16458.93771163884,18455.00933318336,20286.95447549651,22128.17874970573,23977.63207769134,25835.95529949772,27709.96306801069,28259.78882524174,28834.24767053889,29432.91672586709,30035.94630329657,30638.37367670823,31060.22317899247,31485.1492212366,31912.54597126568,32344.95251399587,32802.12458293604,33303.80969695441,33853.97308969918,34468.93363153014,35164.02827533264,35905.30901417151,36660.37924488024,37424.4158851918,38192.55113445189,38960.31013860155,39739.5663142659,40520.9093810313,41501.63893899975,42438.59684048261,43314.00934782792,44163.74219723378,44991.65339536441,45817.74241292285,46634.3157185211,47439.95082215829,48241.71671170378,49022.67157925863,49819.00318622762,50645.04344532066,51514.81435205101,52437.97193142539,53400.32741025716,54417.39429604196,55471.4115232818,56543.69026756018,57625.26060495171,58699.59467562589,59779.07842339089,60863.24377347192,61949.81602896668,63029.55017020732,64071.66164803806,65084.30913486562,66081.71031117698,67076.99793920563,68071.95977009565,69050.23609396606,70026.48610896801,71005.75240086297,71996.32494263421,73009.32511045043,74041.33336871921,75118.12858343754,76244.98633217163,77421.14402575247,78639.47737185529,79866.36707977265,81111.23038064732,82366.23330156185,83627.80858915552,84894.42145535658)
for ($c = 4; $c -le 93; $c++) { $ws2.Cells.Item(2, $c).Value = $row2_2[$c - 4] }

# Sheet 2, row 3
$row2_3 = @(16.48537993724342,35.64126850001205,67.95227137196412,91.55882505369019,154.1114174301326,189.5380254539318,229.0984327837957,274.6454455175687,310.6976431225825,355.1824150807658,376.8012459508082,401.1198250528995,451.3731157110839,502.0669407484772,553.210121837977,683.9536430731865,814.8899566572799,946.144222431652,1077.868034183833,1210.080284410197,1331.22809995027,1452.955683065615,1575.19257351267,1697.981822019387,1821.781050976558,1856.627749878317,1893.084281645027,1931.122194037628,1969.40837794095,2007.609496817973,2033.64112540239,2059.832892577169,2086.143746239999,2112.743690507408,2140.95884482782,2172.11448580981,2206.476639822199,2245.140536992881,2289.135577612189,2336.183291524237,2384.10782881589,2432.58572411339,2481.290563799489,2529.92234418108,2579.276275779511,2628.718300649148,2691.500114530676,2751.293517608481,2806.905871656777,2860.744475941925,2913.062725560857,2965.203245435645,3016.649340745571,3067.305273777461,3117.645122192956,3166.52675978077,3216.37943340978,3268.165001615569,3322.824590630447,3381.00622011332,3441.753710947987,3506.108746787329,3572.879233281455,3640.810855154261,3709.301656288735,3777.237179731276,3845.449243544042,3913.9066173466,3982.456854119564,4050.479700014814,4115.904570836551,4179.279557284064,4241.558947605265,4303.624092062675,4365.594833216182,4426.368212042011,4486.928361722128,4547.613594362457,4608.980191593895,4671.774539662394,4735.761727563869,4802.672391826973,4872.860581099542,4946.275042246549,5022.437542448236,5099.086387510393,5176.854184485947,5255.215155411203,5333.929926017497,5412.895567229443)
for ($c = 4; $c -le 93; $c++) { $ws2.Cells.Item(3, $c).Value = $row2_3[$c - 4] }

# Sheet 2, row 4
$row2_4 = @(804639.2050553659,1353362.529320378,1708436.427766096,1999982.367653959,2327988.870166786,2584097.520784012,2833697.580378841,3071280.737027797,3377249.990599348,3657753.052128214,3964329.110268774,4351078.669680058,6152608.998734247,7961206.851362932,9774067.730484076,15463285.08776202,21154893.70088381,26850065.5192398,32549174.10306183,38252214.73638891,43627102.9057575,49006978.30586433,54390737.266729,59778235.84872077,65170513.14576576,66068168.06564754,66972619.30749371,67885322.80019078,68807068.39199966,69739577.43116887,69890446.40953721,70062634.03627315,70264576.48270591,70507701.14445975,70806279.64061813,71172124.25519237,71611664.30828971,72139423.4210823,72767119.06782754,73466789.31124635,74208833.3030888,74982493.74582787,75775254.55503045,76573692.55590217,77810982.33601242,79156360.11339197,80636535.38390981,82044213.09263304,83357980.22609037,84603306.44176149,85722103.60563241,86801684.33614601,87843653.5036065,88854122.44951001,89848399.65979671,90760396.35430177,91689408.54196812,92650629.35551733,93657315.57835229,94720134.8375407,95783329.10039249,96908200.14266607,98081361.87141028,99287355.88325062,100516991.0952663,101720729.7971047,102937295.2009112,104161369.4962721,105385341.4366284,106595628.8213768,107729574.9029123,108821390.3902808,109880225.7750834,110916913.7759729,111934852.8361636,112889208.6223745,113827697.6347889,114757299.5838688,115687650.9326703,116629616.0990172,117543019.4490939,118487396.05236,119468833.9146859,120489188.5424356,121545080.4851215,122559077.8672679,123591459.5783014,124635872.54712,125688448.0508006,126746095.2632313)
for ($c = 4; $c -le 93; $c++) { $ws2.Cells.Item(4, $c).Value = $row2_4[$c - 4] }

# Sheet 2, row 5
$row2_5 = @(1694765.500019309,2690689.44184167,3189769.617742783,3496332.863593473,3783547.520168905,3986926.552072899,4152788.813036441,4366579.503820671,4755296.09295396,5093711.330577493,5518897.662736825,6052115.534921395,9582036.944963485,13114454.27040061,16649925.46994091,27903211.01515893,39160785.48616605,50423257.10180578,61691193.89669359,72965100.80613247,83737515.32484347,94516718.45160845,105303090.2230472,116097119.9160368,126899578.160813,128700832.9205105,130513927.9883276,132342132.3055398,134190184.8280987,136064529.3629172,136268903.7029659,136517696.6862194,136822097.2611616,137193983.5016896,137645099.0583053,138186060.554462,138825255.9869242,139567735.515134,140414175.7716143,141360055.3992901,142395242.9886493,143504074.4381345,144666042.7886414,145857090.9264696,148200962.485128,150847798.1095374,153449078.0129937,155985591.633616,158443542.8103154,160815352.8021049,162932901.4315465,164968396.041324,166932644.748394,168841039.8639432,170712134.9110026,172401210.416401,174092863.2455202,175805205.1709886,177553552.856917,179349828.5363047,181078475.0435126,182868264.9161348,184720790.1569129,186634640.5859428,188605700.2031895,190541077.3927115,192518742.7262243,194528744.9456218,196560143.2456701,198601562.7935913,200557466.3825096,202501751.4073018,204425386.2204912,206321309.5426199,208184901.6331586,209925773.9184367,211633437.494013,213311552.2406553,214966058.3410127,216604586.2325904,218135931.8404877,219668835.4325197,221211807.0530727,222772409.6281063,224356831.0633768,225833153.6442609,227340547.0563503,228879886.254897,230450298.0409352,232049343.715721)
for ($c = 4; $c -le 93; $c++) { $ws2.Cells.Item(5, $c).Value = $row2_5[$c - 4] }
